$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "On s" + "tudying the performance..." -> "On studying the
#    performance..." (merges the two runs into one and removes the
#    "_GoBack" bookmark that used to sit between them).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "On studying the performance of Hadoop Map Reduce vs MPI for Aggregation Operations: A Big Data Challenge",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "On studying the performance of Hadoop Map Reduce vs MPI for Aggregation Operations: A Big Data Challenge",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "...and compare them. The parameters for improvement would be wall
#    clock time." -> "...and compare them using wall clock time as a
#    performance measure."
#    The target keeps the replacement text split across four runs, so we
#    bracket each piece with temporary bookmarks while editing (bookmarks
#    stop the engine from silently re-coalescing same-formatted runs) and
#    remove the markers again afterwards.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    ". The parameters for improvement would be wall clock time.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $rng.Start
$e = $rng.End

$d.Bookmarks.Add("zzTmp1", $d.Range($s, $s)) | Out-Null
$d.Bookmarks.Add("zzTmp2", $d.Range($s + 2, $s + 2)) | Out-Null
$d.Bookmarks.Add("zzTmp3", $d.Range($s + 42, $s + 42)) | Out-Null
$d.Bookmarks.Add("zzTmp4", $d.Range($s + 57, $s + 57)) | Out-Null
$d.Bookmarks.Add("zzTmp5", $d.Range($e, $e)) | Out-Null

$d.Range($s + 57, $e).Text = "."
$d.Range($s + 42, $s + 57).Text = " as a performance measure"
$d.Range($s + 2, $s + 42).Text = " wall clock time"
$d.Range($s, $s + 2).Text = " using"

$d.Bookmarks("zzTmp1").Delete()
$d.Bookmarks("zzTmp2").Delete()
$d.Bookmarks("zzTmp3").Delete()
$d.Bookmarks("zzTmp4").Delete()
$d.Bookmarks("zzTmp5").Delete()

# ---------------------------------------------------------------------------
# 3) Drop "safety against race conditions, " from the list of challenges.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "in getting the right libraries, debugging network problems, safety against race conditions, and implementation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in getting the right libraries, debugging network problems, and implementation",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Split "...handling library issues (majorly...) and " into "...handling
#    lib" / "rary issues (majorly...) and " with the "_GoBack" bookmark
#    re-inserted at the split point (matching where it now lives after the
#    title was cleaned up above).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    " handling library issues (majorly library version problems or corrupted library builds) and ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $rng2.Start + [string]" handling lib".Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint)) | Out-Null

# ---------------------------------------------------------------------------
# 5) Merge ". However, later due " + "to consideration...Java 8. " back
#    into a single run (also drops the lastRenderedPageBreak that used to
#    sit on the second half).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ". However, later due to consideration of issues of the same package for running the software, Java 9 was degraded to Java 8. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". However, later due to consideration of issues of the same package for running the software, Java 9 was degraded to Java 8. ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Merge "decided...ingest the data. " + "Hence, data would be ingested
#    by using " into one run (drops the lastRenderedPageBreak that used to
#    sit on the second half).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "decided that a layer of abstraction by using another software is not a good way to ingest the data. Hence, data would be ingested by using ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "decided that a layer of abstraction by using another software is not a good way to ingest the data. Hence, data would be ingested by using ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Refresh the "Dbengines" run so it drops its stale lastRenderedPageBreak
#    marker (text is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dbengines", $true, $false, $false, $false, $false, $true, 1, $false,
    "Dbengines", 2) | Out-Null
